$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-24 Tuesday" "2024-09-25 Wednesday"

Replace-Text "595÷9=" "217÷3="
Replace-Text "679÷5=" "762÷2="
Replace-Text "336÷2=" "207÷7="
Replace-Text "997÷2=" "405÷7="
Replace-Text "618÷4=" "776÷9="
Replace-Text "662÷5=" "862÷8="
Replace-Text "331÷4=" "461÷4="
Replace-Text "586÷6=" "437÷2="
Replace-Text "439÷6=" "335÷8="
Replace-Text "869÷7=" "310÷9="
Replace-Text "457÷2=" "735÷9="
Replace-Text "428÷8=" "598÷3="
Replace-Text "143÷8=" "598÷9="
Replace-Text "409÷5=" "218÷5="
Replace-Text "595÷6=" "843÷8="
Replace-Text "794÷3=" "296÷6="
Replace-Text "671÷8=" "178÷3="
Replace-Text "591÷3=" "810÷5="
Replace-Text "210÷6=" "144÷2="
Replace-Text "371÷5=" "437÷5="
Replace-Text "957÷9=" "873÷3="
Replace-Text "244÷2=" "345÷9="
Replace-Text "261÷2=" "967÷9="
Replace-Text "126÷8=" "441÷3="
Replace-Text "666÷5=" "180÷7="
